$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Add the new column K (year 2020) with values copied/derived from column J --
$ws.Range("K3").Value = 2020
$ws.Range("K4").Value = 0
$ws.Range("K5").Value = 48.2
$ws.Range("K6").Value = 19.3
$ws.Range("K7").Value = 24.2
$ws.Range("K8").Value = 8.3

# -- Update the active cell / selection shown in the workbook --
$ws.Range("J22").Select()
